# Auto-generated edit script: update KHL probabilities tour data for 2025-12-16 matches
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Probabilities")

# Row 2
$ws.Range("A2").Value = 1369
$ws.Range("B2").Value = '2025-12-16T15:30:00'
$ws.Range("C2").Value = 'Сибирь'
$ws.Range("D2").Value = 'Авангард'
$ws.Range("E2").Value = 897869
$ws.Range("F2").Value = 'https://text.khl.ru/text/897869.html'
$ws.Range("G2").Value = 1.066667
$ws.Range("H2").Value = 5.6
$ws.Range("I2").Value = 3.845946
$ws.Range("J2").Value = 2.015022
$ws.Range("K2").Value = 1.540844
$ws.Range("L2").Value = 4.722973
$ws.Range("M2").Value = 6.666667
$ws.Range("N2").Value = 22.624
$ws.Range("O2").Value = 39.66778
$ws.Range("P2").Value = 62.29178
$ws.Range("Q2").Value = -0.2
$ws.Range("R2").Value = 0.2
$ws.Range("S2").Value = 0.061991
$ws.Range("T2").Value = 0.073958
$ws.Range("U2").Value = 0.85472
$ws.Range("V2").Value = 0.129169
$ws.Range("W2").Value = 0.861501
$ws.Range("X2").Value = 0.251294
$ws.Range("Y2").Value = 0.739376
$ws.Range("Z2").Value = 0.404287
$ws.Range("AA2").Value = 0.586383
$ws.Range("AB2").Value = 0.564008
$ws.Range("AC2").Value = 0.426662
$ws.Range("AD2").Value = 0.706931
$ws.Range("AE2").Value = 0.28374
$ws.Range("AF2").Value = 0.455751
$ws.Range("AG2").Value = 0.544249
$ws.Range("AH2").Value = 0.201474
$ws.Range("AI2").Value = 0.798526
$ws.Range("AJ2").Value = 0.94913
$ws.Range("AK2").Value = 0.05087
$ws.Range("AL2").Value = 0.849992
$ws.Range("AM2").Value = 0.150008
$ws.Range("AN2").Value = 0.25272
$ws.Range("AO2").Value = 0.966774

# Row 3
$ws.Range("A3").Value = 1369
$ws.Range("B3").Value = '2025-12-16T17:00:00'
$ws.Range("C3").Value = 'Автомобилист'
$ws.Range("D3").Value = 'Барыс'
$ws.Range("E3").Value = 897868
$ws.Range("F3").Value = 'https://text.khl.ru/text/897868.html'
$ws.Range("G3").Value = 5.3
$ws.Range("H3").Value = 1.035714
$ws.Range("I3").Value = 3.197554
$ws.Range("J3").Value = 5.035714
$ws.Range("K3").Value = 5.167857
$ws.Range("L3").Value = 2.116634
$ws.Range("M3").Value = 6.335714
$ws.Range("N3").Value = 35.22868
$ws.Range("O3").Value = 24.47065
$ws.Range("P3").Value = 59.699331
$ws.Range("Q3").Value = 0.2
$ws.Range("R3").Value = -0.2
$ws.Range("S3").Value = 0.814377
$ws.Range("T3").Value = 0.081043
$ws.Range("U3").Value = 0.087563
$ws.Range("V3").Value = 0.068088
$ws.Range("W3").Value = 0.914895
$ws.Range("X3").Value = 0.148584
$ws.Range("Y3").Value = 0.834399
$ws.Range("Z3").Value = 0.265858
$ws.Range("AA3").Value = 0.717126
$ws.Range("AB3").Value = 0.408237
$ws.Range("AC3").Value = 0.574746
$ws.Range("AD3").Value = 0.556404
$ws.Range("AE3").Value = 0.426579
$ws.Range("AF3").Value = 0.964863
$ws.Range("AG3").Value = 0.035137
$ws.Range("AH3").Value = 0.888792
$ws.Range("AI3").Value = 0.111208
$ws.Range("AJ3").Value = 0.624644
$ws.Range("AK3").Value = 0.375356
$ws.Range("AL3").Value = 0.354858
$ws.Range("AM3").Value = 0.645142
$ws.Range("AN3").Value = 0.943187
$ws.Range("AO3").Value = 0.28523

# Row 4
$ws.Range("A4").Value = 1369
$ws.Range("B4").Value = '2025-12-16T18:00:00'
$ws.Range("C4").Value = 'Лада'
$ws.Range("D4").Value = 'ЦСКА'
$ws.Range("E4").Value = 897871
$ws.Range("F4").Value = 'https://text.khl.ru/text/897871.html'
$ws.Range("G4").Value = 0.933333
$ws.Range("H4").Value = 1.743933
$ws.Range("I4").Value = 1.979672
$ws.Range("J4").Value = 1.361111
$ws.Range("K4").Value = 1.147222
$ws.Range("L4").Value = 1.861803
$ws.Range("M4").Value = 2.677266
$ws.Range("N4").Value = 20.472803
$ws.Range("O4").Value = 28.115807
$ws.Range("P4").Value = 48.58861
$ws.Range("Q4").Value = -0.2
$ws.Range("R4").Value = -0.093748
$ws.Range("S4").Value = 0.232994
$ws.Range("T4").Value = 0.226297
$ws.Range("U4").Value = 0.540705
$ws.Range("V4").Value = 0.64521
$ws.Range("W4").Value = 0.354786
$ws.Range("X4").Value = 0.813745
$ws.Range("Y4").Value = 0.186251
$ws.Range("Z4").Value = 0.915169
$ws.Range("AA4").Value = 0.084826
$ws.Range("AB4").Value = 0.966034
$ws.Range("AC4").Value = 0.033961
$ws.Range("AD4").Value = 0.987899
$ws.Range("AE4").Value = 0.012096
$ws.Range("AF4").Value = 0.318219
$ws.Range("AG4").Value = 0.681781
$ws.Range("AH4").Value = 0.109274
$ws.Range("AI4").Value = 0.890726
$ws.Range("AJ4").Value = 0.555298
$ws.Range("AK4").Value = 0.444702
$ws.Range("AL4").Value = 0.285979
$ws.Range("AM4").Value = 0.714021
$ws.Range("AN4").Value = 0.691107
$ws.Range("AO4").Value = 0.909845

# Row 5
$ws.Range("A5").Value = 1369
$ws.Range("B5").Value = '2025-12-16T19:00:00'
$ws.Range("C5").Value = 'Локомотив'
$ws.Range("D5").Value = 'Динамо Мн'
$ws.Range("E5").Value = 897866
$ws.Range("F5").Value = 'https://text.khl.ru/text/897866.html'
$ws.Range("G5").Value = 2.233075
$ws.Range("H5").Value = 4.376873
$ws.Range("I5").Value = 1.525518
$ws.Range("J5").Value = 1.214286
$ws.Range("K5").Value = 1.72368
$ws.Range("L5").Value = 2.951196
$ws.Range("M5").Value = 6.609949
$ws.Range("N5").Value = 25.39464
$ws.Range("O5").Value = 35.477968
$ws.Range("P5").Value = 60.872609
$ws.Range("Q5").Value = -0.128841
$ws.Range("R5").Value = 0.170011
$ws.Range("S5").Value = 0.206709
$ws.Range("T5").Value = 0.164605
$ws.Range("U5").Value = 0.628429
$ws.Range("V5").Value = 0.313656
$ws.Range("W5").Value = 0.686087
$ws.Range("X5").Value = 0.499263
$ws.Range("Y5").Value = 0.500479
$ws.Range("Z5").Value = 0.672802
$ws.Range("AA5").Value = 0.326941
$ws.Range("AB5").Value = 0.808014
$ws.Range("AC5").Value = 0.191729
$ws.Range("AD5").Value = 0.898314
$ws.Range("AE5").Value = 0.101429
$ws.Range("AF5").Value = 0.514073
$ws.Range("AG5").Value = 0.485927
$ws.Range("AH5").Value = 0.249041
$ws.Range("AI5").Value = 0.750959
$ws.Range("AJ5").Value = 0.793443
$ws.Range("AK5").Value = 0.206557
$ws.Range("AL5").Value = 0.565787
$ws.Range("AM5").Value = 0.434213
$ws.Range("AN5").Value = 0.560992
$ws.Range("AO5").Value = 0.903817

# Row 6
$ws.Range("A6").Value = 1369
$ws.Range("B6").Value = '2025-12-16T19:00:00'
$ws.Range("C6").Value = 'Ак Барс'
$ws.Range("D6").Value = 'Салават Юлаев'
$ws.Range("E6").Value = 897867
$ws.Range("F6").Value = 'https://text.khl.ru/text/897867.html'
$ws.Range("G6").Value = 3.926991
$ws.Range("H6").Value = 1.05
$ws.Range("I6").Value = 3.197954
$ws.Range("J6").Value = 2.252974
$ws.Range("K6").Value = 3.089982
$ws.Range("L6").Value = 2.123977
$ws.Range("M6").Value = 4.976991
$ws.Range("N6").Value = 33.449071
$ws.Range("O6").Value = 22.215265
$ws.Range("P6").Value = 55.664336
$ws.Range("Q6").Value = 0.081373
$ws.Range("R6").Value = -0.2
$ws.Range("S6").Value = 0.577475
$ws.Range("T6").Value = 0.16553
$ws.Range("U6").Value = 0.256607
$ws.Range("V6").Value = 0.236266
$ws.Range("W6").Value = 0.763346
$ws.Range("X6").Value = 0.403786
$ws.Range("Y6").Value = 0.595827
$ws.Range("Z6").Value = 0.578474
$ws.Range("AA6").Value = 0.421138
$ws.Range("AB6").Value = 0.730277
$ws.Range("AC6").Value = 0.269336
$ws.Range("AD6").Value = 0.843347
$ws.Range("AE6").Value = 0.156265
$ws.Range("AF6").Value = 0.813895
$ws.Range("AG6").Value = 0.186105
$ws.Range("AH6").Value = 0.596665
$ws.Range("AI6").Value = 0.403335
$ws.Range("AJ6").Value = 0.626512
$ws.Range("AK6").Value = 0.373488
$ws.Range("AL6").Value = 0.356839
$ws.Range("AM6").Value = 0.643161
$ws.Range("AN6").Value = 0.865995
$ws.Range("AO6").Value = 0.601063

# Row 7
$ws.Range("A7").Value = 1369
$ws.Range("B7").Value = '2025-12-16T19:00:00'
$ws.Range("C7").Value = 'Нефтехимик'
$ws.Range("D7").Value = 'Торпедо'
$ws.Range("E7").Value = 897872
$ws.Range("F7").Value = 'https://text.khl.ru/text/897872.html'
$ws.Range("G7").Value = 3.691483
$ws.Range("H7").Value = 2.343097
$ws.Range("I7").Value = 4.704545
$ws.Range("J7").Value = 2.371485
$ws.Range("K7").Value = 3.031484
$ws.Range("L7").Value = 3.523821
$ws.Range("M7").Value = 6.03458
$ws.Range("N7").Value = 35.326163
$ws.Range("O7").Value = 31.500113
$ws.Range("P7").Value = 66.826275
$ws.Range("Q7").Value = 0.169913
$ws.Range("R7").Value = -0.06921
$ws.Range("S7").Value = 0.346666
$ws.Range("T7").Value = 0.156411
$ws.Range("U7").Value = 0.49553
$ws.Range("V7").Value = 0.1081
$ws.Range("W7").Value = 0.890506
$ws.Range("X7").Value = 0.217553
$ws.Range("Y7").Value = 0.781053
$ws.Range("Z7").Value = 0.361053
$ws.Range("AA7").Value = 0.637554
$ws.Range("AB7").Value = 0.517833
$ws.Range("AC7").Value = 0.480773
$ws.Range("AD7").Value = 0.664654
$ws.Range("AE7").Value = 0.333953
$ws.Range("AF7").Value = 0.805505
$ws.Range("AG7").Value = 0.194495
$ws.Range("AH7").Value = 0.583827
$ws.Range("AI7").Value = 0.416173
$ws.Range("AJ7").Value = 0.866608
$ws.Range("AK7").Value = 0.133392
$ws.Range("AL7").Value = 0.683536
$ws.Range("AM7").Value = 0.316464
$ws.Range("AN7").Value = 0.658213
$ws.Range("AO7").Value = 0.785403

# Row 8
$ws.Range("A8").Value = 1369
$ws.Range("B8").Value = '2025-12-16T19:30:00'
$ws.Range("C8").Value = 'Динамо М'
$ws.Range("D8").Value = 'Спартак'
$ws.Range("E8").Value = 897865
$ws.Range("F8").Value = 'https://text.khl.ru/text/897865.html'
$ws.Range("G8").Value = 2.579904
$ws.Range("H8").Value = 5.1
$ws.Range("I8").Value = 1.21875
$ws.Range("J8").Value = 6.2
$ws.Range("K8").Value = 4.389952
$ws.Range("L8").Value = 3.159375
$ws.Range("M8").Value = 7.679904
$ws.Range("N8").Value = 29.315697
$ws.Range("O8").Value = 36.075524
$ws.Range("P8").Value = 65.391221
$ws.Range("Q8").Value = -0.024741
$ws.Range("R8").Value = 0.2
$ws.Range("S8").Value = 0.597759
$ws.Range("T8").Value = 0.134552
$ws.Range("U8").Value = 0.261647
$ws.Range("V8").Value = 0.057255
$ws.Range("W8").Value = 0.936703
$ws.Range("X8").Value = 0.128506
$ws.Range("Y8").Value = 0.865452
$ws.Range("Z8").Value = 0.236086
$ws.Range("AA8").Value = 0.757873
$ws.Range("AB8").Value = 0.371444
$ws.Range("AC8").Value = 0.622514
$ws.Range("AD8").Value = 0.517425
$ws.Range("AE8").Value = 0.476533
$ws.Range("AF8").Value = 0.933157
$ws.Range("AG8").Value = 0.066843
$ws.Range("AH8").Value = 0.81366
$ws.Range("AI8").Value = 0.18634
$ws.Range("AJ8").Value = 0.823425
$ws.Range("AK8").Value = 0.176575
$ws.Range("AL8").Value = 0.611553
$ws.Range("AM8").Value = 0.388447
$ws.Range("AN8").Value = 0.83849
$ws.Range("AO8").Value = 0.543736

# Row 9
$ws.Range("A9").Value = 1369
$ws.Range("B9").Value = '2025-12-16T19:30:00'
$ws.Range("C9").Value = 'СКА'
$ws.Range("D9").Value = 'Драконы'
$ws.Range("E9").Value = 897870
$ws.Range("F9").Value = 'https://text.khl.ru/text/897870.html'
$ws.Range("G9").Value = 2.979392
$ws.Range("H9").Value = 2.257346
$ws.Range("I9").Value = 2.143498
$ws.Range("J9").Value = 5.526316
$ws.Range("K9").Value = 4.252854
$ws.Range("L9").Value = 2.200422
$ws.Range("M9").Value = 5.236738
$ws.Range("N9").Value = 32.012917
$ws.Range("O9").Value = 29.012809
$ws.Range("P9").Value = 61.025726
$ws.Range("Q9").Value = -0.018638
$ws.Range("R9").Value = -0.049881
$ws.Range("S9").Value = 0.725148
$ws.Range("T9").Value = 0.117991
$ws.Range("U9").Value = 0.152386
$ws.Range("V9").Value = 0.115106
$ws.Range("W9").Value = 0.880419
$ws.Range("X9").Value = 0.228944
$ws.Range("Y9").Value = 0.766581
$ws.Range("Z9").Value = 0.375869
$ws.Range("AA9").Value = 0.619655
$ws.Range("AB9").Value = 0.533895
$ws.Range("AC9").Value = 0.46163
$ws.Range("AD9").Value = 0.679578
$ws.Range("AE9").Value = 0.315947
$ws.Range("AF9").Value = 0.925286
$ws.Range("AG9").Value = 0.074714
$ws.Range("AH9").Value = 0.796656
$ws.Range("AI9").Value = 0.203344
$ws.Range("AJ9").Value = 0.645533
$ws.Range("AK9").Value = 0.354467
$ws.Range("AL9").Value = 0.377399
$ws.Range("AM9").Value = 0.622601
$ws.Range("AN9").Value = 0.920724
$ws.Range("AO9").Value = 0.42033

